$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the date column (A) with the same "previous day + 1" shared formula
# pattern used for rows 4-10, continuing it down through row 28.
$ws.Range("A11:A28").Formula = "=A10+1"

# Fill in the URL column (B) for the newly added rows. The writes are ordered
# so that each distinct URL is first introduced in the same sequence as the
# shared-string table in the target workbook (rows 14/15 and 20/21 have their
# string table entries swapped relative to row order).
$ws.Range("B11").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/1ab5a678c34cce2fa4e095f158cc1b152bac64d0/counties.json"
$ws.Range("B12").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/7daa747efd2b02f3e1bd61c2d9844fe7761bb88d/counties.json"
$ws.Range("B13").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/59dc3c257afd9d129702f0fa3cac73945923eb75/counties.json"
$ws.Range("B15").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/8f182bd89dfc7efac6033cdff502dead207a5c9a/counties.json"
$ws.Range("B14").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/2b88b4a7afeb95cc81c4e4bf834a3a8981c94285/counties.json"
$ws.Range("B16").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/b46a9e4137a23a86b08770d3c737ee9ae84a051a/counties.json"
$ws.Range("B17").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/4d055f8550e1c5635d7425632b7be582edee6377/counties.json"
$ws.Range("B18").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/06ed24aa6756de13a6d1d0c283469c9d5c0d25d1/counties.json"
$ws.Range("B19").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/8482fd1f83a924b4da62b43afd96f1c7d3828ffe/counties.json"
$ws.Range("B21").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/ed6993ff4256dbe6b347dcd69c9151e2b60185fb/counties.json"
$ws.Range("B20").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/272fb781e1fcb50d0b7e32a22d39a6a0e84d8756/counties.json"
$ws.Range("B22").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/05a863ecc79a7568e2cae96170ba08ad957ba885/counties.json"
$ws.Range("B23").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/8b5ae6c27bf498701eb3b33b217714f9f2f9d8c0/counties.json"
$ws.Range("B24").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/93a4e1f47993b9b81c1a5851dbc6839e3f4707f8/counties.json"
$ws.Range("B25").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/c73948dd87ed84471e50eeb13e92efe255b943cf/counties.json"
$ws.Range("B26").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/f194e99b69bbc39ae8387b236c9041aa442f6bb1/counties.json"
$ws.Range("B27").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/0ba1dd7c6c6eadf91f897c48286751b7b5f2b297/counties.json"
$ws.Range("B28").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/1d1e010216c917e56d8fb9f072fb10e3a129f2f7/counties.json"

# Match the saved selection state from the source workbook.
$ws.Range("O22").Select()
